$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.021.86'
$ws.Range("E2").Value = '  +1.43%  '

# Row 3
$ws.Range("D3").Value = '1.851.66'
$ws.Range("E3").Value = '  +1.37%  '

# Row 5
$ws.Range("D5").Value = '''1.012'
$ws.Range("E5").Value = '  +0.44%  '

# Row 6
$ws.Range("D6").Value = '''310.11'
$ws.Range("E6").Value = '  +0.47%  '

# Row 7
$ws.Range("D7").Value = '''0.4779'
$ws.Range("E7").Value = '  +2.15%  '

# Row 8
$ws.Range("E8").Value = '  +2.09%  '

# Row 9
$ws.Range("D9").Value = '''0.07226'
$ws.Range("E9").Value = '  +1.31%  '

# Row 10
$ws.Range("D10").Value = '''0.9328'
$ws.Range("E10").Value = '  +3.29%  '

# Row 11
$ws.Range("D11").Value = '''19.75'
$ws.Range("E11").Value = '  +1.68%  '

# Row 12
$ws.Range("D12").Value = '''0.07727'
$ws.Range("E12").Value = '  -0.33%  '

# Row 13
$ws.Range("D13").Value = '1.865.81'
$ws.Range("E13").Value = '  +2.05%  '

# Row 14
$ws.Range("D14").Value = '''5.329'
$ws.Range("E14").Value = '  +1.00%  '

# Row 15
$ws.Range("E15").Value = '  +1.19%  '

# Row 16
$ws.Range("D16").Value = '''88.99'
$ws.Range("E16").Value = '  +1.66%  '

# Row 17
$ws.Range("D17").Value = '''1.015'
$ws.Range("E17").Value = '  +0.55%  '

# Row 18
$ws.Range("D18").Value = '''0.000008645'
$ws.Range("E18").Value = '  +1.11%  '

# Row 19
$ws.Range("D19").Value = '''1.012'
$ws.Range("E19").Value = '  +0.46%  '

# Row 20
$ws.Range("D20").Value = '27.065.59'
$ws.Range("E20").Value = '  +1.46%  '

# Row 21
$ws.Range("D21").Value = '''14.55'
$ws.Range("E21").Value = '  +2.26%  '

# Row 22
$ws.Range("D22").Value = '''5.058'
$ws.Range("E22").Value = '  +0.65%  '

# Row 23
$ws.Range("D23").Value = '''10.67'
$ws.Range("E23").Value = '  +1.09%  '

# Row 24
$ws.Range("E24").Value = '  +1.58%  '

# Row 25
$ws.Range("D25").Value = '''152.80'
$ws.Range("E25").Value = '  -0.11%  '

# Row 26
$ws.Range("E26").Value = '  +1.69%  '

# Row 27
$ws.Range("D27").Value = '''2.009'
$ws.Range("E27").Value = '  +1.50%  '

# Row 28
$ws.Range("D28").Value = '''114.53'
$ws.Range("E28").Value = '  +0.52%  '

# Row 29
$ws.Range("D29").Value = '''4.996'
$ws.Range("E29").Value = '  +2.62%  '

# Row 30
$ws.Range("D30").Value = '''0.08904'
$ws.Range("E30").Value = '  +1.04%  '

# Row 31
$ws.Range("D31").Value = '''3.319'
$ws.Range("E31").Value = '  +5.49%  '

# Row 32
$ws.Range("E32").Value = '  +0.94%  '

# Row 33
$ws.Range("D33").Value = '''0.7458'
$ws.Range("E33").Value = '  +1.18%  '

# Row 34
$ws.Range("D34").Value = '''4.508'
$ws.Range("E34").Value = '  +1.41%  '

# Row 35
$ws.Range("D35").Value = '''2.735'
$ws.Range("E35").Value = '  -3.93%  '

# Row 36
$ws.Range("D36").Value = '''1.110'
$ws.Range("E36").Value = '  +2.76%  '

# Row 37
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").Value = '''0.01956'
$ws.Range("E37").Value = '  +1.23%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.05278'
$ws.Range("E38").Value = '  +2.26%  '

# Row 39
$ws.Range("D39").Value = '''2.976'
$ws.Range("E39").Value = '  +2.60%  '

# Row 40
$ws.Range("D40").Value = '''0.5227'
$ws.Range("E40").Value = '  +3.36%  '

# Row 41
$ws.Range("D41").Value = '''7.024'
$ws.Range("E41").Value = '  +2.08%  '

# Row 42
$ws.Range("D42").Value = '''0.1514'
$ws.Range("E42").Value = '  +1.17%  '

# Row 43
$ws.Range("D43").Value = '''8.230'
$ws.Range("E43").Value = '  +2.29%  '

# Row 44
$ws.Range("D44").Value = '''10.58'
$ws.Range("E44").Value = '  +6.01%  '

# Row 45
$ws.Range("D45").Value = '''0.4751'
$ws.Range("E45").Value = '  +1.79%  '

# Row 46
$ws.Range("D46").Value = '''1.014'
$ws.Range("E46").Value = '  +0.53%  '

# Row 47
$ws.Range("D47").Value = '''101.64'
$ws.Range("E47").Value = '  +3.70%  '

# Row 48
$ws.Range("D48").Value = '''1.612'
$ws.Range("E48").Value = '  +2.42%  '

# Row 49
$ws.Range("D49").Value = '''65.70'
$ws.Range("E49").Value = '  +2.60%  '

# Row 50
$ws.Range("D50").Value = '''0.06035'
$ws.Range("E50").Value = '  -0.24%  '

# Row 51
$ws.Range("D51").Value = '''0.8891'
$ws.Range("E51").Value = '  +4.21%  '

